$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for "Macroferia Regional de Talca" /
# Acelga. It belongs at the top of the data block (row 169), so push the
# existing rows 169-214 down to 170-215 and populate the freed row with the
# new observation.
$ws.Rows(169).Insert()

$ws.Cells.Item(169, 1).Value = 5
$ws.Cells.Item(169, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(169, 3).Value = 'Maule'
$ws.Cells.Item(169, 4).Value = 44642
$ws.Cells.Item(169, 5).Value = 7
$ws.Cells.Item(169, 6).Value = 100112009
$ws.Cells.Item(169, 7).Value = 'Acelga'
$ws.Cells.Item(169, 8).Value = 'Sin especificar'
$ws.Cells.Item(169, 9).Value = 'Primera'
$ws.Cells.Item(169, 10).Value = 400
$ws.Cells.Item(169, 11).Value = 3500
$ws.Cells.Item(169, 12).Value = 3500
$ws.Cells.Item(169, 13).Value = 3500
$ws.Cells.Item(169, 14).Value = '$/docena de atados (4 kilos)'
$ws.Cells.Item(169, 15).Value = 'Región del Maule'
$ws.Cells.Item(169, 16).Value = 875
$ws.Cells.Item(169, 17).Value = 4
$ws.Cells.Item(169, 18).Value = 'Hortaliza'
